# Cashflow model: replace the formula-driven value in J6 with a fixed
# (hard-coded) constant. J6 previously mirrored I5/J5 via a shared
# formula ("=J5", evaluating to "-"); the model now pins this to a fixed
# cashflow figure. Every later row in column J (J7:J29) just copies the
# cell above it (=J6, =J7, ...), so they automatically pick up the new
# constant on recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J6").Value = 100000000

# Column J now holds a wide number instead of the short "-" placeholder,
# so widen it to fit the new content.
$ws.Columns("J:J").AutoFit() | Out-Null

# Leave the selection where the edit was reviewed.
$ws.Range("J7").Select() | Out-Null
